$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.230.48'
$ws.Range("E2").Value = '  +14.84%  '

$ws.Range("D3").Value = '1.675.77'
$ws.Range("E3").Value = '  +8.91%  '

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = '  -0.69%  '

$ws.Range("D5").Value = "'307.55"
$ws.Range("E5").Value = '  +9.05%  '

$ws.Range("D6").Value = "'0.9986"
$ws.Range("E6").Value = '  +3.20%  '

$ws.Range("E7").Value = '  +2.64%  '

$ws.Range("D8").Value = "'0.3425"
$ws.Range("E8").Value = '  +7.57%  '

$ws.Range("D9").Value = "'47.76"
$ws.Range("E9").Value = '  +17.46%  '

$ws.Range("D10").Value = "'1.176"
$ws.Range("E10").Value = '  +7.16%  '

$ws.Range("D11").Value = "'0.07271"
$ws.Range("E11").Value = '  +6.66%  '

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = '  -0.32%  '

$ws.Range("D13").Value = "'20.44"
$ws.Range("E13").Value = '  +8.84%  '

$ws.Range("D14").Value = "'6.080"
$ws.Range("E14").Value = '  +7.08%  '

$ws.Range("D15").Value = "'6.742"
$ws.Range("E15").Value = '  +5.93%  '

$ws.Range("D16").Value = '1.679.62'
$ws.Range("E16").Value = '  +9.37%  '

$ws.Range("D17").Value = "'0.00001103"
$ws.Range("E17").Value = '  +5.80%  '

$ws.Range("D18").Value = "'0.9990"
$ws.Range("E18").Value = '  +3.24%  '

$ws.Range("D19").Value = "'0.06701"
$ws.Range("E19").Value = '  +9.73%  '

$ws.Range("D20").Value = "'81.28"
$ws.Range("E20").Value = '  +12.28%  '

$ws.Range("D21").Value = "'16.37"
$ws.Range("E21").Value = '  +8.53%  '

$ws.Range("E22").Value = '  +7.32%  '

$ws.Range("D23").Value = "'12.01"
$ws.Range("E23").Value = '  +5.77%  '

$ws.Range("D24").Value = '24.240.05'
$ws.Range("E24").Value = '  +14.57%  '

$ws.Range("D25").Value = "'2.405"
$ws.Range("E25").Value = '  +3.69%  '

$ws.Range("D26").Value = "'3.360"
$ws.Range("E26").Value = '  -9.12%  '

$ws.Range("D27").Value = "'2.649"
$ws.Range("E27").Value = '  +19.07%  '

$ws.Range("D28").Value = "'151.93"
$ws.Range("E28").Value = '  +2.26%  '

$ws.Range("D29").Value = "'19.45"
$ws.Range("E29").Value = '  +10.11%  '

$ws.Range("D30").Value = '1.863.31'
$ws.Range("E30").Value = '  +9.30%  '

$ws.Range("D31").Value = "'126.68"
$ws.Range("E31").Value = '  +6.98%  '

$ws.Range("D32").Value = "'6.411"
$ws.Range("E32").Value = '  +23.19%  '

$ws.Range("D33").Value = "'4.051"
$ws.Range("E33").Value = '  +0.94%  '

$ws.Range("D34").Value = "'0.9829"
$ws.Range("E34").Value = '  +15.05%  '

$ws.Range("D35").Value = "'1.744"
$ws.Range("E35").Value = '  +15.22%  '

$ws.Range("D36").Value = "'0.08450"
$ws.Range("E36").Value = '  +5.56%  '

$ws.Range("E37").Value = '  +16.89%  '

$ws.Range("D38").Value = "'5.361"
$ws.Range("E38").Value = '  +8.78%  '

$ws.Range("D39").Value = "'0.06402"
$ws.Range("E39").Value = '  +9.43%  '

$ws.Range("D40").Value = "'8.806"
$ws.Range("E40").Value = '  +14.18%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = "'0.02332"
$ws.Range("E41").Value = '  +10.74%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = "'1.279"
$ws.Range("E42").Value = '  +5.70%  '

$ws.Range("D43").Value = "'0.2097"
$ws.Range("E43").Value = '  +9.19%  '

$ws.Range("D44").Value = "'0.6135"
$ws.Range("E44").Value = '  +12.65%  '

$ws.Range("D45").Value = "'0.9968"
$ws.Range("E45").Value = '  +3.09%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'13.13"
$ws.Range("E46").Value = '  +4.71%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = "'3.798"
$ws.Range("E47").Value = '  +6.43%  '

$ws.Range("D48").Value = "'0.5933"
$ws.Range("E48").Value = '  +9.06%  '

$ws.Range("D49").Value = "'127.95"
$ws.Range("E49").Value = '  +4.82%  '

$ws.Range("D50").Value = "'2.011"
$ws.Range("E50").Value = '  +7.57%  '

$ws.Range("D51").Value = "'0.07167"
$ws.Range("E51").Value = '  +9.20%  '
